$p = $ppt.ActivePresentation
try {
    $d = $p.Designs.Add("theme2")
    Write-Output "Add OK: $d, Count=$($p.Designs.Count)"
} catch {
    Write-Output "Add ERR: $($_.Exception.Message)"
}
